# dplyr Assignment Instructions - sync syllabus: bump homework number,
# tighten the "System Environment" line and note the sysInfo() command.

$d = $word.ActiveDocument

# 1. Title: "Homework 4:  dplyr" -> "Homework 5:  dplyr"
$d.Content.Find.Execute("Homework 4", $true, $false, $false, $false, $false, $true, 1, $false, "Homework 5", 2)

# 2. "System Environment" line: drop the extra space before "(1/2 point)"
#    (keeping the existing italics on "(1/2 point)" intact) ...
$d.Content.Find.Execute("Environment  ", $true, $false, $false, $false, $false, $true, 1, $false, "Environment ", 2)
# ... and append a note about the sysInfo() command after the trailing
# spaces (which stay in their own, non-italic run).
$d.Content.Find.Execute(")    ", $true, $false, $false, $false, $false, $true, 1, $false, ")    - use sysInfo() command.", 2)
